# Adjust bulk upload files: add a "Tipo" column (H) classifying each
# sensor certificate row as Primario / Secundario / Vencido.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("H1").Value = "Tipo"

# First populate every data row with the common "Secundario" value so
# that shared string gets registered right after "Tipo" ...
for ($row = 2; $row -le 22; $row++) {
    $ws.Cells.Item($row, 8).Value = "Secundario"
}

# ... then apply the "Vencido" override (row 21) before any "Primario"
# override, matching the order new strings were first introduced in the
# saved workbook (Tipo, Secundario, Vencido, Primario).
$ws.Cells.Item(21, 8).Value = "Vencido"

# Finally set the "Primario" rows.
$ws.Cells.Item(6, 8).Value = "Primario"
$ws.Cells.Item(9, 8).Value = "Primario"
$ws.Cells.Item(22, 8).Value = "Primario"

# Match the cursor/selection state recorded in the saved workbook.
$ws.Range("J6").Select()
